# STV-Project simple_test_data.xlsx : "ignore merge cells (View Title)"
#
# View1 ("Main" section) gets a merged, centered title row reading "Main"
# inserted above the tab_list/folder_list rows, and a second merged,
# centered title row reading "Edit" inserted below them (followed by one
# blank spacer row). View2 stops being the active/selected sheet; View1
# becomes active instead, and the remembered selections on both sheets
# change to reflect where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "View1"
$ws2 = $wb.Worksheets.Item(2)   # "View2"

# --- View1: make room for the two new title rows -------------------------
# Original layout:
#   row1 headers, row2 tab_list, row3 folder_list
# Target layout:
#   row1 headers, row2 "Main" (new), row3 tab_list, row4 folder_list,
#   row5 "Edit" (new), row6 blank spacer (new)
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(5).Insert()
$ws1.Rows.Item(6).Insert()

# "Main" title row, merged & centered across A2:C2
$ws1.Range("A2:C2").ClearFormats()
$ws1.Range("A2:C2").HorizontalAlignment = -4108
$ws1.Range("A2").Value = "Main"
$ws1.Range("A2:C2").Merge()

# "Edit" title row, merged & centered across A5:C5
$ws1.Range("A5:C5").ClearFormats()
$ws1.Range("A5:C5").HorizontalAlignment = -4108
$ws1.Range("A5").Value = "Edit"
$ws1.Range("A5:C5").Merge()

# Blank spacer cell left behind on row 6
$ws1.Range("B6").ClearFormats()

# --- View2: remember the new selection, but it's no longer the active tab
$ws2.Range("B41").Select()

# --- View1 becomes the active sheet / tab, with its own new selection ----
$ws1.Activate()
$ws1.Range("A5:C5").Select()
